# "Generate Report for Handoff" -- the localization status report has moved
# from "In Translation" to "Ready for handoff"; refresh the status cells and
# the associated handoff timestamps on all three sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet --------------------------------------------------------
# E2/F2 hold the per-language status ("zh-cn" / "de-de" columns), G2 holds the
# "Latest HO Xliff Generate Date" timestamp.
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-15 10:39:11"

# --- zh-cn detail sheet ------------------------------------------------------
# C2 = Status, H2 = Latest Handoff Datetime
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-15 10:39:05"

# --- de-de detail sheet ------------------------------------------------------
# C2 = Status, H2 = Latest Handoff Datetime
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-15 10:39:11"

# --- Column widths -----------------------------------------------------------
# "Ready for handoff" is wider than "In Translation", so the report generator
# widened the Status columns to fit the new text.
$overview.Columns.Item(5).ColumnWidth = 16.38
$overview.Columns.Item(6).ColumnWidth = 16.38
$zhcn.Columns.Item(3).ColumnWidth = 16.38
$dede.Columns.Item(3).ColumnWidth = 16.38
